$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Cornucopia"
$ws.Range("C2").Value = "None"

$ws.Range("B3").Value = "Cornucopia Guilds"
$ws.Range("C3").Value = "Update Pack"

$ws.Range("B4").Value = "Dark Ages"
$ws.Range("C4").Value = "None"

$ws.Range("B5").Value = "Nocturne"
$ws.Range("C5").Value = "None"

$ws.Range("B6").Value = "Rising Sun"
$ws.Range("C6").Value = "None"
